$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, F (number), G (intercept), I (sse), T (predict)
$data = @(
    @(2, 5125415.8, 1.4949104, 4.8216911, 22853906.22374598),
    @(3, 8308127.3, 1.4949104, 4.8216911, 37045416.33659144),
    @(4, 6532697.6, 1.4949104, 4.8216911, 29128887.13116513),
    @(5, 1471564.9, 1.4949104, 4.8216911, 6561615.201396174),
    @(6, 2214795.8, 1.4949104, 4.8216911, 9875634.971497621),
    @(7, 5066818.3, 1.4949104, 4.8216911, 22592623.66205685),
    @(8, 5678551.2, 1.4949104, 4.8216911, 25320302.13266604),
    @(9, 3774931.2, 1.4949104, 4.8216911, 16832180.45370932),
    @(10, 4393006.3, 1.4949104, 4.8216911, 19588138.39464991),
    @(11, 2955146.5, 1.4949104, 4.8216911, 13176812.11121079),
    @(12, 1946882.9, 1.4949104, 4.8216911, 8681028.225108024),
    @(13, 5868404.8, 1.4949104, 4.8216911, 26166847.37698369),
    @(14, 7185126.5, 1.4949104, 4.8216911, 32038026.50250389),
    @(15, 10212215, 1.2033945, 7.8131369, 34021036.464082),
    @(16, 17565000, 1.2033945, 7.8131369, 58516150.07044017),
    @(17, 16974378, 1.2033945, 7.8131369, 56548548.27215361),
    @(18, 9932519.6, 1.2033945, 7.8131369, 33089257.47174429),
    @(19, 5696901.4, 1.2033945, 7.8131369, 18978692.69905498),
    @(20, 8308127.3, 1.2033945, 7.8131369, 27677746.87680734),
    @(21, 10465713, 1.2033945, 7.8131369, 34865541.27538609),
    @(22, 5462777.8, 1.2033945, 7.8131369, 18198731.88421685),
    @(23, 3122602.5, 1.2033945, 7.8131369, 10402657.35840203),
    @(24, 3138636.4, 1.2033945, 7.8131369, 10456072.79242506),
    @(25, 8118576.2, 1.2033945, 7.8131369, 27046275.16524361),
    @(26, 8973280.5, 1.2033945, 7.8131369, 29893642.37757782),
    @(27, 5161862.9, 1.2033945, 7.8131369, 17196262.11781597),
    @(28, 5365959.4, 1.2033945, 7.8131369, 17876190.46525984),
    @(29, 4609687.2, 1.2033945, 7.8131369, 15356740.56208296),
    @(30, 4020323.6, 1.2033945, 7.8131369, 13393331.00537047),
    @(31, 5469063.3, 1.2033945, 7.8131369, 18219671.43794687),
    @(32, 7075820.9, 1.2033945, 7.8131369, 23572433.68380057),
    @(33, 7079015.1, 1.2033945, 7.8131369, 23583074.86151506),
    @(34, 5386637.3, 1.2033945, 7.8131369, 17945076.94971994),
    @(35, 8040267.7, 1.2033945, 7.8131369, 26785397.74208442)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 6).Value = $row[1]
    $ws.Cells.Item($r, 7).Value = $row[2]
    $ws.Cells.Item($r, 9).Value = $row[3]
    $ws.Cells.Item($r, 20).Value = $row[4]
}

